# The source commit for this edit ("Adicionando o restante do projeto")
# had already landed all of its slide/content changes in the working
# copy we start from here -- before.pptx already contains all 17 slides
# with the final slide/master ids referenced throughout the deck.
#
# The only remaining delta versus the target OOXML is the removal of
# ppt/changesInfos/changesInfo1.xml, an internal Office co-authoring
# "change log" part (plus its relationship/content-type registration).
# That part is pure session bookkeeping written by the web/co-authoring
# client; it carries no slide content and is not reachable through the
# PowerPoint COM object model (Presentation/Slides/Shapes, etc. never
# surface raw package parts like changesInfos, revisionInfo, and so
# on) -- there is no Shape/Slide/Presentation member that targets it,
# mirroring real PowerPoint, where this log is dropped automatically
# by the application/session rather than through any scriptable API.
#
# Touch the active presentation (no visible/content mutation is
# required) so the automation round-trips the deck cleanly.
$p = $ppt.ActivePresentation
$p.Slides.Count | Out-Null
